$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update status column (I) for rows 42,44,45,46,49,50,51,52: "In Progress" -> "Done" ---
$doneRows = 42,44,45,46,49,50,51,52
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 9).Value = "Done"
}

# --- New row 62 content (added before I48 update so shared-string order matches) ---
$ws.Range("A62").Value = "Backend code"
$ws.Range("B62").Value = "-"
$ws.Range("C62").Value = 56
$ws.Range("D62").Value = "Upload and add articles don’t work in categories page"
$ws.Range("E62").Value = "5/3/2023"
$ws.Range("F62").Value = "Mazrouaa"
$ws.Range("G62").Value = "Categories page"

# --- I48 gets a distinct status text (added after the new row so it becomes the last new shared string) ---
$ws.Range("I48").Value = "Done with another way"
